$d = $word.ActiveDocument

# The document currently has this paragraph layout (1-indexed):
#   1  : "I.A. de Ogre" (title)
#   2-16: assorted "Sensores/Ataques/Bateria" notes + blank paragraphs
#   17 : "Objetivo"
#   18 : "    A la hora de desarrollar..." (long paragraph)
#   19 : "    Por otra parte, ..." (long paragraph, ends with the _GoBack bookmark)
#   20+: unchanged remainder of the document
#
# The edit removes paragraphs 2-16 entirely (the old class-notes block) and
# replaces them with a single empty paragraph that now carries the
# _GoBack bookmark (moved from the end of paragraph 19).

# 1) Delete the whole block of paragraphs between the title and "Objetivo".
$delStart = $d.Paragraphs.Item(2).Range.Start
$delEnd = $d.Paragraphs.Item(16).Range.End
$d.Range($delStart, $delEnd).Delete()

# 2) Insert a paragraph break right before "Objetivo" (now paragraph 2),
#    producing a new, empty paragraph between the title and "Objetivo".
$objStart = $d.Paragraphs.Item(2).Range.Start
$d.Range($objStart, $objStart).InsertParagraphBefore()

# 3) Re-home the _GoBack bookmark inside that new empty paragraph.
#    Insert a temporary placeholder character so the collapsed bookmark
#    range is unambiguously anchored inside the new paragraph (and not at
#    its boundary with a neighbour), bookmark around it, then remove the
#    placeholder text again, leaving a clean, run-less paragraph.
$newParaStart = $d.Paragraphs.Item(2).Range.Start
$d.Range($newParaStart, $newParaStart).InsertAfter("X")
$bmRange = $d.Range($newParaStart, $newParaStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($newParaStart, $newParaStart + 1).Delete()
